{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright block (and the blank\n// paragraph immediately preceding it) that follows the last requirement\n// paragraph (\"LOQ4233: Gest\u00e3o de Neg\u00f3cios (Requisito fraco)\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the first target paragraph so we can also remove the\n// blank paragraph immediately before it.\nlet firstTargetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    firstTargetIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (firstTargetIndex > 0 && items[firstTargetIndex - 1].text === \"\") {\n  toDelete.push(items[firstTargetIndex - 1]);\n}\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.includes(items[i].text)) {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright block (and the blank\n# paragraph immediately preceding it) that follows the last requirement\n# paragraph (\"LOQ4233: Gest\u00e3o de Neg\u00f3cios (Requisito fraco)\").\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd([char]13, [char]7)\n}\n\n$markerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n$idx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq $markerText) {\n        $idx = $i\n        break\n    }\n}\n\nif ($idx -gt 0) {\n    # Delete the blank paragraph right before the marker paragraph, if present.\n    if ($idx -gt 1 -and (Get-ParaText $d.Paragraphs.Item($idx - 1)) -eq \"\") {\n        $d.Paragraphs.Item($idx - 1).Range.Delete()\n        $idx = $idx - 1\n    }\n\n    # The marker paragraph (\"Ver no Jupiter...\") is now at $idx again\n    # (indices shifted down by one after the delete above only if we deleted\n    # the preceding paragraph; recompute defensively).\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ((Get-ParaText $d.Paragraphs.Item($i)) -eq $markerText) {\n            $idx = $i\n            break\n        }\n    }\n\n    # Delete the \"Ver no Jupiter...\" paragraph and the copyright paragraph\n    # that follows it. Deleting the same index twice removes both, since\n    # each deletion shifts the following paragraph into that slot.\n    $d.Paragraphs.Item($idx).Range.Delete()\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
